# Update "想去人数" (F) and "最低票价" (G) values on the "展览" and "全部类型"
# sheets, which hold identical data. Only F2/G2 change both columns; all
# other touched rows only change column F.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value  = 2070
    $ws.Range("G2").Value  = 65

    $ws.Range("F7").Value  = 1710

    $ws.Range("F9").Value  = 705

    $ws.Range("F11").Value = 106

    $ws.Range("F14").Value = 231

    $ws.Range("F19").Value = 3965

    $ws.Range("F24").Value = 940

    $ws.Range("F25").Value = 644

    $ws.Range("F27").Value = 10

    $ws.Range("F29").Value = 1768

    $ws.Range("F30").Value = 30

    $ws.Range("F32").Value = 60

    $ws.Range("F33").Value = 180
}
